# Mar 9 - Content Updates
# Insert a new 5-row block ("Prod / Smileactives / Core2" + "End") at the
# top of the data table on the "Content" sheet, pushing all existing
# blocks down by 5 rows (row 2 -> row 7, etc.), matching the layout already
# used for the other blocks (2 content rows followed by 3 blank rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Content")

# Make room: push everything from row 2 down by 5 rows.
$ws.Rows("2:6").Insert()

# New block content.
$ws.Range("A2").Value = "Prod"
$ws.Range("B2").Value = "Smileactives"
$ws.Range("C2").Value = "Core2"
$ws.Range("A3").Value = "End"

# Match the shaded formatting used by the other data/end rows (copy the
# format from the equivalent rows that follow, which already use it).
$ws.Range("A7:C7").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Rows 4-6 (and B3:C3) stay completely empty - no stray styled cells.
$ws.Range("B3:C6").Clear()
$ws.Range("A4:A6").Clear()

# Match the author's final selection state (row 3 selected).
$null = $ws.Rows("3:3").Select()

Write-Output "done"
